# Adds a "graph" adjacency list puzzle (and the preceding "arrange all 0s to
# left" puzzle heading) to the Learning sheet, and removes the now-unused
# Call/Put option-pricing scratch table that lived on Sheet6.

$wb = $excel.ActiveWorkbook

# --- Learning sheet: new puzzle content -----------------------------------
$ws = $wb.Worksheets.Item("Learning")

$ws.Range("A34").Value = "Arrange all 0 to left in the given array"

$ws.Range("A43").Value = "Graph"

$ws.Range("B44").Value = "A"
$ws.Range("C44").Formula = '="->B,C"'

$ws.Range("B45").Value = "B"
$ws.Range("C45").Formula = '="->A,C,D"'

$ws.Range("B46").Value = "C"
$ws.Range("C46").Formula = '="->B,D"'

$ws.Range("B47").Value = "D"
$ws.Range("C47").Formula = '="->B,C"'

$ws.Range("B48").Value = "E"
$ws.Range("C48").Formula = '="->F"'

$ws.Range("B49").Value = "F"
$ws.Range("C49").Formula = '="->E"'

$ws.Activate()
$ws.Range("I39").Select()

# --- Trip sheet: just a view-state nudge (selection moved) ----------------
$ws2 = $wb.Worksheets.Item("Trip")
$ws2.Activate()
$ws2.Range("C31").Select()

# --- Sheet6: the Call/Put options table is no longer needed, clear it -----
$ws4 = $wb.Worksheets.Item("Sheet6")
$ws4.Cells.Clear()
$ws4.Activate()
$ws4.Range("F36").Select()

# Re-activate Learning so it stays the selected tab, matching the source
$ws.Activate()
